$d = $word.ActiveDocument

# Locate the old sub-domain text ("dohmh") inside the contact e-mail and
# replace it with the new one ("nyc"), e.g.
#   rrohrer@health.dohmh.gov  ->  rrohrer@health.nyc.gov
$rng = $d.Content
$found = $rng.Find.Execute("dohmh", $true, $false, $false, $false, $false, $true, 1, $false, "nyc", 2)

# The engine auto-merges newly-written text back into the surrounding run
# when it shares identical formatting. Force the replacement text into its
# own run (matching the source edit, which split the single run into three)
# by nudging a character attribute away from, then back to, the shared
# formatting -- the two transient writes bracket the insertion in its own
# run without leaving any visible formatting difference behind.
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("nyc", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng2.Bold = 1
$rng2.Bold = 0
